$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Id=9, "Decision Tree Regression") gets an updated accuracy value.
$ws.Range("C11").Value = 81.049431756350501

# Insert a new row 12 for "Decision Tree Regression+Feature Selection",
# copying the border/format of the row above it (row 11).
$ws.Rows.Item(12).Insert()
$ws.Range("A11:C11").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)
[void]($excel.CutCopyMode = 0)

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Decision Tree Regression+Feature Selection"
$ws.Range("C12").Value = 81.739564185612196

# Update the active selection to match the edited workbook's saved state.
[void]$ws.Range("E7").Select()

Write-Output "done"
